# Update Efna3-Epha7 NATMI LR-pair sheet with new TPM-based values.
# "Neutrophils" sending-cluster rows (7-11) are now "MuSCs", and all
# recomputed ligand/receptor expression + specificity statistics are refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("M2").Value = 0.2371515
$ws.Range("N2").Value = 0.474303
$ws.Range("O2").Value = 0.1229805122625808
$ws.Range("P2").Value = 0.1147457658281868
$ws.Range("Q2").Value = 0.04283442250574999
$ws.Range("R2").Value = 0.171337690023
$ws.Range("S2").Value = 0.08958526814769611
$ws.Range("T2").Value = 0.08358665947481668

# Row 3
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("O3").Value = 0.05352595923371758
$ws.Range("P3").Value = 0.07491281022046142
$ws.Range("Q3").Value = 0.01864322656216666
$ws.Range("R3").Value = 0.111859359373
$ws.Range("S3").Value = 0.03899103461674435
$ws.Range("T3").Value = 0.05457030604140266

# Row 4
$ws.Range("G4").Value = 0.1806205
$ws.Range("H4").Value = 0.361241
$ws.Range("I4").Value = 0.7284509268949775
$ws.Range("J4").Value = 0.7284509268949775
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04845666666666667
$ws.Range("N4").Value = 0.14537
$ws.Range("O4").Value = 0.02512834913211087
$ws.Range("P4").Value = 0.03516864109744933
$ws.Range("Q4").Value = 0.008752267361666666
$ws.Range("R4").Value = 0.05251360417
$ws.Range("S4").Value = 0.01830476921662676
$ws.Range("T4").Value = 0.02561862920507376

# Row 5
$ws.Range("G5").Value = 0.1806205
$ws.Range("H5").Value = 0.361241
$ws.Range("I5").Value = 0.7284509268949775
$ws.Range("J5").Value = 0.7284509268949775
$ws.Range("M5").Value = 1.414436
$ws.Range("N5").Value = 2.828872
$ws.Range("O5").Value = 0.7334891992782493
$ws.Range("P5").Value = 0.6843749334706176
$ws.Range("Q5").Value = 0.255476137538
$ws.Range("R5").Value = 1.021904550152
$ws.Range("S5").Value = 0.5343108870816955
$ws.Range("T5").Value = 0.49853355463036

# Row 6
$ws.Range("G6").Value = 0.1806205
$ws.Range("H6").Value = 0.361241
$ws.Range("I6").Value = 0.7284509268949775
$ws.Range("J6").Value = 0.7284509268949775
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1251046666666667
$ws.Range("N6").Value = 0.375314
$ws.Range("O6").Value = 0.06487598009334153
$ws.Range("P6").Value = 0.09079784938328471
$ws.Range("Q6").Value = 0.02259646744566667
$ws.Range("R6").Value = 0.135578804674
$ws.Range("S6").Value = 0.04725896783221475
$ws.Range("T6").Value = 0.0661417775433243

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.067331
$ws.Range("H7").Value = 0.134662
$ws.Range("I7").Value = 0.2715490731050226
$ws.Range("J7").Value = 0.2715490731050226
$ws.Range("M7").Value = 0.2371515
$ws.Range("N7").Value = 0.474303
$ws.Range("O7").Value = 0.1229805122625808
$ws.Range("P7").Value = 0.1147457658281868
$ws.Range("Q7").Value = 0.0159676476465
$ws.Range("R7").Value = 0.063870590586
$ws.Range("S7").Value = 0.03339524411488468
$ws.Range("T7").Value = 0.03115910635337009

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.067331
$ws.Range("H8").Value = 0.134662
$ws.Range("I8").Value = 0.2715490731050226
$ws.Range("J8").Value = 0.2715490731050226
$ws.Range("O8").Value = 0.05352595923371758
$ws.Range("P8").Value = 0.07491281022046142
$ws.Range("Q8").Value = 0.006949748714333332
$ws.Range("R8").Value = 0.041698492286
$ws.Range("S8").Value = 0.01453492461697323
$ws.Range("T8").Value = 0.02034250417905876

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.067331
$ws.Range("H9").Value = 0.134662
$ws.Range("I9").Value = 0.2715490731050226
$ws.Range("J9").Value = 0.2715490731050226
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04845666666666667
$ws.Range("N9").Value = 0.14537
$ws.Range("O9").Value = 0.02512834913211087
$ws.Range("P9").Value = 0.03516864109744933
$ws.Range("Q9").Value = 0.003262635823333334
$ws.Range("R9").Value = 0.01957581494
$ws.Range("S9").Value = 0.006823579915484105
$ws.Range("T9").Value = 0.00955001189237557

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.067331
$ws.Range("H10").Value = 0.134662
$ws.Range("I10").Value = 0.2715490731050226
$ws.Range("J10").Value = 0.2715490731050226
$ws.Range("M10").Value = 1.414436
$ws.Range("N10").Value = 2.828872
$ws.Range("O10").Value = 0.7334891992782493
$ws.Range("P10").Value = 0.6843749334706176
$ws.Range("Q10").Value = 0.09523539031600001
$ws.Range("R10").Value = 0.380941561264
$ws.Range("S10").Value = 0.1991783121965538
$ws.Range("T10").Value = 0.1858413788402577

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.067331
$ws.Range("H11").Value = 0.134662
$ws.Range("I11").Value = 0.2715490731050226
$ws.Range("J11").Value = 0.2715490731050226
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1251046666666667
$ws.Range("N11").Value = 0.375314
$ws.Range("O11").Value = 0.06487598009334153
$ws.Range("P11").Value = 0.09079784938328471
$ws.Range("Q11").Value = 0.008423422311333334
$ws.Range("R11").Value = 0.050540533868
$ws.Range("S11").Value = 0.01761701226112679
$ws.Range("T11").Value = 0.02465607183996041
